$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $r = $cell.Range
    # Trim the trailing cell-end mark(s) so we only overwrite the visible text.
    $r.MoveEnd(1, -1) | Out-Null
    $r.Text = $newText
}

Set-CellText $t 6 3 "cognition"
Set-CellText $t 6 5 "behavior"
Set-CellText $t 9 4 "heart rate"
Set-CellText $t 10 3 "stimulus"
Set-CellText $t 10 4 "visual"
Set-CellText $t 10 6 "adult"
Set-CellText $t 10 7 "visual"
Set-CellText $t 11 3 "human"
Set-CellText $t 11 4 "scene"
Set-CellText $t 11 6 "electrophysiological"
Set-CellText $t 11 7 "control"
